# Add Lec 15 CI reading link, renumber Lec 17/18 rows after removing a
# duplicate "Lec 17:" row, and add new "Lec 19: " / "Lec 20:" placeholder
# rows, matching the commit "add Lec 15 CI".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: add the "15-CI" reading link next to "Lec 15: ..." ---
# Copy the formatting of the neighboring already-styled D26 cell so the new
# E26 cell picks up the same (non-bold Helvetica) style used throughout the
# table, then set its text.
$ws.Range("D26").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "15-CI"

# --- Row 28: remove the stray "Lec 17:" label (it moves down to row 29) ---
$ws.Range("C28").Clear() | Out-Null

# --- Rows 29/30: shift the lecture labels up by one slot ---
$ws.Range("C29").Value = "Lec 17:"
$ws.Range("C30").Value = "Lec 18:"

# --- Rows 32/33: fill in the two new lecture labels ---
$ws.Range("C32").Value = "Lec 19: "
$ws.Range("C33").Value = "Lec 20:"

# --- Update the window scroll position / active selection ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C33").Select() | Out-Null
